$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1741.5172
$ws.Range("J17").Value = 1913.8182
$ws.Range("L17").Value = 5741.4546
$ws.Range("N17").Value = -6077.4546

$ws.Range("H33").Value = 306.4
$ws.Range("I33").Value = 133
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 133
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 96
$ws.Range("N33").Value = -1458

$ws.Range("H40").Value = 4944.1665
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350

$ws.Range("H62").Value = 16124.68
$ws.Range("I62").Value = 15209.526
$ws.Range("J62").Value = 19022.666
$ws.Range("K62").Value = 15209.526
$ws.Range("L62").Value = 19022.666
$ws.Range("M62").Value = -14585.526
$ws.Range("N62").Value = -20270.666

$ws.Range("H65").Value = 16124.68
$ws.Range("I65").Value = 15209.526
$ws.Range("J65").Value = 19022.666
$ws.Range("K65").Value = 76047.63
$ws.Range("L65").Value = 95113.33
$ws.Range("M65").Value = -72927.63
$ws.Range("N65").Value = -101353.33

$ws.Range("H76").Value = 4451.4546
$ws.Range("I76").Value = 4441.1113
$ws.Range("K76").Value = 4441.1113
$ws.Range("M76").Value = -4126.1113

$ws.Range("H79").Value = 4451.4546
$ws.Range("I79").Value = 4441.1113
$ws.Range("K79").Value = 4441.1113
$ws.Range("M79").Value = -3349.1113

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H87").Value = 95150.5
$ws.Range("J87").Value = 95150.5
$ws.Range("L87").Value = 95150.5
$ws.Range("N87").Value = -97646.5

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H90").Value = 95150.5
$ws.Range("J90").Value = 95150.5
$ws.Range("L90").Value = 285451.5
$ws.Range("N90").Value = -297931.5

$ws.Range("H94").Value = 3892.9285
$ws.Range("I94").Value = 2269.3076
$ws.Range("J94").Value = 25000
$ws.Range("K94").Value = 2269.3076
$ws.Range("L94").Value = 25000
$ws.Range("M94").Value = -1818.3076
$ws.Range("N94").Value = -25902

$ws.Range("H100").Value = 4388.5557
$ws.Range("I100").Value = 1232.2941
$ws.Range("J100").Value = 9754.2
$ws.Range("K100").Value = 1232.2941
$ws.Range("L100").Value = 9754.2
$ws.Range("M100").Value = -691.2941000000001
$ws.Range("N100").Value = -10836.2

$ws.Range("H132").Value = 361094.28
$ws.Range("I132").Value = 375618.9
$ws.Range("K132").Value = 1126856.7
$ws.Range("M132").Value = -1124326.7

$ws.Range("H135").Value = 27028242
$ws.Range("I135").Value = 29412912
$ws.Range("J135").Value = 1966.3334
$ws.Range("K135").Value = 264716208
$ws.Range("L135").Value = 17697.0006
$ws.Range("M135").Value = -264713673
$ws.Range("N135").Value = -22767.0006

$ws.Range("H137").Value = 3296.3333
$ws.Range("I137").Value = 1719.8
$ws.Range("K137").Value = 5159.4
$ws.Range("M137").Value = -2609.4

$ws.Range("H138").Value = 3414.0833
$ws.Range("J138").Value = 3601.4355
$ws.Range("L138").Value = 10804.3065
$ws.Range("N138").Value = -21084.3065

$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9268747
$ws.Range("I32").Value = 10875343
$ws.Range("K32").Value = 10875343
$ws.Range("M32").Value = -10875056

$ws.Range("H37").Value = 73199.3
$ws.Range("J37").Value = 78999.75
$ws.Range("L37").Value = 78999.75
$ws.Range("N37").Value = -79545.75

$ws.Range("H43").Value = 42166.5
$ws.Range("I43").Value = 43113.668
$ws.Range("J43").Value = 41598.2
$ws.Range("K43").Value = 43113.668
$ws.Range("L43").Value = 41598.2
$ws.Range("M43").Value = -42800.668
$ws.Range("N43").Value = -42224.2

$ws.Range("H44").Value = 15069500
$ws.Range("J44").Value = 15069500
$ws.Range("L44").Value = 15069500
$ws.Range("N44").Value = -15070476

$ws.Range("H55").Value = 149996
$ws.Range("J55").Value = 149996
$ws.Range("L55").Value = 149996
$ws.Range("N55").Value = -150626

$ws.Range("H61").Value = 7674.643
$ws.Range("I61").Value = 3530.0625
$ws.Range("K61").Value = 3530.0625
$ws.Range("M61").Value = -3318.0625

$ws.Range("H74").Value = 4035055
$ws.Range("I74").Value = 5436625
$ws.Range("K74").Value = 5436625
$ws.Range("M74").Value = -5435751

$ws.Range("H77").Value = 4035055
$ws.Range("I77").Value = 5436625
$ws.Range("K77").Value = 27183125
$ws.Range("M77").Value = -27178757

$ws.Range("H80").Value = 127947.5
$ws.Range("J80").Value = 127947.5
$ws.Range("L80").Value = 127947.5
$ws.Range("N80").Value = -129943.5

$ws.Range("H83").Value = 127947.5
$ws.Range("J83").Value = 127947.5
$ws.Range("L83").Value = 383842.5
$ws.Range("N83").Value = -393826.5

$ws.Range("H132").Value = 434274.2
$ws.Range("I132").Value = 685903.06
$ws.Range("K132").Value = 2057709.18
$ws.Range("M132").Value = -2055179.18

$ws.Range("H136").Value = 7674.643
$ws.Range("I136").Value = 3530.0625
$ws.Range("K136").Value = 10590.1875
$ws.Range("M136").Value = -8040.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2087.318
$ws.Range("I20").Value = 2615.3333
$ws.Range("J20").Value = 1453.7
$ws.Range("K20").Value = 2615.3333
$ws.Range("L20").Value = 1453.7
$ws.Range("M20").Value = -2368.3333
$ws.Range("N20").Value = -1947.7

$ws.Range("H82").Value = 33946.5
$ws.Range("I82").Value = 9990.6
$ws.Range("J82").Value = 63891.375
$ws.Range("K82").Value = 9990.6
$ws.Range("L82").Value = 63891.375
$ws.Range("M82").Value = -9607.6
$ws.Range("N82").Value = -64657.375

$ws.Range("H85").Value = 33946.5
$ws.Range("I85").Value = 9990.6
$ws.Range("J85").Value = 63891.375
$ws.Range("K85").Value = 9990.6
$ws.Range("L85").Value = 63891.375
$ws.Range("M85").Value = -8664.6
$ws.Range("N85").Value = -66543.375

$ws.Range("H86").Value = 2085.2222
$ws.Range("I86").Value = 2085.2222
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2085.2222
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -962.2222000000002
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 2085.2222
$ws.Range("I89").Value = 2085.2222
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10426.111
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4810.111000000001
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 1307.909
$ws.Range("I94").Value = 1462.5714
$ws.Range("J94").Value = 1037.25
$ws.Range("K94").Value = 1462.5714
$ws.Range("L94").Value = 1037.25
$ws.Range("M94").Value = -1011.5714
$ws.Range("N94").Value = -1939.25

$ws.Range("H99").Value = 4159.8
$ws.Range("I99").Value = 3100
$ws.Range("J99").Value = 4866.3335
$ws.Range("K99").Value = 3100
$ws.Range("L99").Value = 4866.3335
$ws.Range("M99").Value = -1602
$ws.Range("N99").Value = -7862.3335

$ws.Range("H134").Value = 774335.4
$ws.Range("I134").Value = 997300.06
$ws.Range("K134").Value = 2991900.18
$ws.Range("M134").Value = -2989365.18

$ws.Range("H138").Value = 87766
$ws.Range("J138").Value = 87766
$ws.Range("L138").Value = 87766
$ws.Range("N138").Value = -98046

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13955.48
$ws.Range("I31").Value = 2499.5
$ws.Range("K31").Value = 2499.5
$ws.Range("M31").Value = -2204.5

$ws.Range("H34").Value = 13955.48
$ws.Range("I34").Value = 2499.5
$ws.Range("K34").Value = 2499.5
$ws.Range("M34").Value = -2297.5

$ws.Range("H50").Value = 119997.5
$ws.Range("J50").Value = 119997.5
$ws.Range("L50").Value = 119997.5
$ws.Range("N50").Value = -121247.5

$ws.Range("H58").Value = 777244.06
$ws.Range("I58").Value = 955069.6
$ws.Range("K58").Value = 955069.6
$ws.Range("M58").Value = -954866.6

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 95765.5
$ws.Range("J60").Value = 124750
$ws.Range("L60").Value = 124750
$ws.Range("N60").Value = -125772

$ws.Range("H68").Value = 84749.5
$ws.Range("J68").Value = 96666
$ws.Range("L68").Value = 96666
$ws.Range("N68").Value = -98164

$ws.Range("H71").Value = 84749.5
$ws.Range("J71").Value = 96666
$ws.Range("L71").Value = 289998
$ws.Range("N71").Value = -297486

$ws.Range("H136").Value = 777244.06
$ws.Range("I136").Value = 955069.6
$ws.Range("K136").Value = 2865208.8
$ws.Range("M136").Value = -2862658.8

$ws.Range("H141").Value = 239206.2
$ws.Range("J141").Value = 287007.88
$ws.Range("L141").Value = 287007.88
$ws.Range("N141").Value = -297367.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 54
$ws.Range("J33").Value = 80
$ws.Range("L33").Value = 480
$ws.Range("N33").Value = -1046

$ws.Range("H92").Value = 680.875
$ws.Range("I92").Value = 763.9286
$ws.Range("J92").Value = 99.5
$ws.Range("K92").Value = 2291.7858
$ws.Range("L92").Value = 298.5
$ws.Range("M92").Value = -1043.7858
$ws.Range("N92").Value = -2794.5

$ws.Range("H97").Value = 512
$ws.Range("I97").Value = 382.66666
$ws.Range("K97").Value = 1147.99998
$ws.Range("M97").Value = -651.99998

$ws.Range("H122").Value = 1048.6757
$ws.Range("J122").Value = 1051.7222
$ws.Range("L122").Value = 9465.4998
$ws.Range("N122").Value = -14365.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9385.417
$ws.Range("I43").Value = 5253.5713
$ws.Range("J43").Value = 15170
$ws.Range("K43").Value = 5253.5713
$ws.Range("L43").Value = 15170
$ws.Range("M43").Value = -5102.5713
$ws.Range("N43").Value = -15472

$ws.Range("H46").Value = 42739.75
$ws.Range("J46").Value = 50306
$ws.Range("L46").Value = 50306
$ws.Range("N46").Value = -50618

$ws.Range("H57").Value = 149997
$ws.Range("J57").Value = 149997
$ws.Range("L57").Value = 149997
$ws.Range("N57").Value = -151637

$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H126").Value = 643802.2
$ws.Range("I126").Value = 982897.94
$ws.Range("K126").Value = 2948693.82
$ws.Range("M126").Value = -2946223.82

$ws.Range("H135").Value = 104996.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 104996.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 104996.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -115136.5

$ws.Range("H140").Value = 104569.71
$ws.Range("J140").Value = 104569.71
$ws.Range("L140").Value = 104569.71
$ws.Range("N140").Value = -114929.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3673.9
$ws.Range("I46").Value = 1332.5
$ws.Range("K46").Value = 1332.5
$ws.Range("M46").Value = -1144.5

$ws.Range("H74").Value = 65000
$ws.Range("I74").Value = 65000
$ws.Range("K74").Value = 65000
$ws.Range("M74").Value = -64002

$ws.Range("H77").Value = 65000
$ws.Range("I77").Value = 65000
$ws.Range("K77").Value = 195000
$ws.Range("M77").Value = -190008

$ws.Range("H100").Value = 7220.9585
$ws.Range("I100").Value = 1608.25
$ws.Range("J100").Value = 12833.667
$ws.Range("K100").Value = 1608.25
$ws.Range("L100").Value = 12833.667
$ws.Range("M100").Value = -1067.25
$ws.Range("N100").Value = -13915.667

$ws.Range("H132").Value = 808126.5
$ws.Range("I132").Value = 1083788.8
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 3251366.4
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -3248836.4
$ws.Range("N132").Value = -23660

$ws.Range("H136").Value = 6628.1816
$ws.Range("I136").Value = 2483.5
$ws.Range("K136").Value = 7450.5
$ws.Range("M136").Value = -4900.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 4
$ws.Range("M14").Value = 164

$ws.Range("H81").Value = 2110
$ws.Range("I81").Value = 2110
$ws.Range("K81").Value = 4220
$ws.Range("M81").Value = -3159

$ws.Range("H84").Value = 2110
$ws.Range("I84").Value = 2110
$ws.Range("K84").Value = 21100
$ws.Range("M84").Value = -15796

$ws.Range("H122").Value = 3506.818
$ws.Range("I122").Value = 3119.5
$ws.Range("K122").Value = 9358.5
$ws.Range("M122").Value = -6908.5

$ws.Range("H126").Value = 3147
$ws.Range("I126").Value = 2829.6667
$ws.Range("J126").Value = 4099
$ws.Range("K126").Value = 8489.000100000001
$ws.Range("L126").Value = 12297
$ws.Range("M126").Value = -6019.000100000001
$ws.Range("N126").Value = -17237

$ws.Range("H135").Value = 85357.5
$ws.Range("J135").Value = 85357.5
$ws.Range("L135").Value = 85357.5
$ws.Range("N135").Value = -95497.5
